$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jun")

# Fill in the newly recorded status rows (135-142) on the "Jun" sheet.
# Column layout: A=#, B=difficulty/category, C=value, D=problem id, E=problem name, F=date

$rows = @(
    @{ Row = 135; B = "시뮬레이션과 구현"; C = 1; D = 19635; E = "배열 돌리기3";   F = 44232 },
    @{ Row = 136; B = "시뮬레이션과 구현"; C = 2; D = 19626; E = "배열 돌리기1";   F = 44234 },
    @{ Row = 137; B = "시뮬레이션과 구현"; C = 2; D = 19627; E = "배열 돌리기2";   F = 44234 },
    @{ Row = 138; B = "시뮬레이션과 구현"; C = 2; D = 14499; E = "주사위 굴리기"; F = 44234 },
    @{ Row = 139; B = "시뮬레이션과 구현"; C = 3; D = 14890; E = "경사로";       F = 44234 },
    @{ Row = 140; B = "시뮬레이션과 구현"; C = 1; D = 15662; E = "톱니바퀴(2)";   F = 44234 },
    @{ Row = 141; B = "시뮬레이션과 구현"; C = 1; D = 14503; E = "로봇 청소기"; F = 44234 },
    @{ Row = 142; B = "시뮬레이션과 구현"; C = 2; D = 15685; E = "드래곤 커브";   F = 44235 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value2 = $r.F
    $ws.Range("F134").Copy() | Out-Null
    $ws.Cells.Item($r.Row, 6).PasteSpecial(-4122)
}

# Restore the view state that Excel recorded after this edit
# (scrolled so row 132 / column C is the top-left visible cell,
# with E137 as the active selection).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 132
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E137").Select()
